$d = $word.ActiveDocument

# Near the end of the document there are three consecutive paragraphs:
#   1. "85, 92, 101, ... 84, 93, 99"    (a list of rating counts)
#   2. an empty paragraph
#   3. "5 stars - ... waste your time"  (a star-rating legend, last paragraph)
#
# The edit collapses all three into a single empty paragraph (the first
# paragraph's formatting survives, but its run/text is removed; the other
# two paragraphs disappear entirely).

$numbersText = "85, 92, 101, 110, 123, 130, 142, 155, 162, 169, 174, 80, 88, 95, 103, 112, 120, 132, 145, 158, 165, 170, 82, 90, 98, 107, 115, 125, 135, 148, 152, 160, 172, 177, 84, 93, 99"

# Find the paragraph holding the ratings-count list by scanning paragraph text
# (robust to any index shifting elsewhere in the document).
$numbersIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*$numbersText*") {
        $numbersIndex = $i
        break
    }
}

if ($numbersIndex -gt 0) {
    # The empty paragraph immediately follows the numbers paragraph, and the
    # "5 stars" paragraph is the document's final paragraph right after that.
    # Remove the empty (non-last) paragraph first -- this host only
    # collapses the document's very last paragraph once per pass, so
    # clearing the earlier one first leaves the star-rating paragraph as the
    # new last paragraph, which can then be removed outright.
    $emptyIndex = $numbersIndex + 1
    $emptyPara = $d.Paragraphs.Item($emptyIndex)
    $emptyPara.Range.Delete()

    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Range.Delete()
}

# Finally, clear the run/text of the numbers paragraph, leaving the empty
# paragraph (with its original formatting) behind.
$d.Content.Find.Execute($numbersText, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
